# ---------------------------------------------------------------------------
# 1) Refresh the panel-query timestamps (column F) on the "data" sheet.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:21:44.711321",
    "2021-10-05 14:21:44.711329",
    "2021-10-05 14:21:44.711332",
    "2021-10-05 14:21:44.711335",
    "2021-10-05 14:21:44.711337",
    "2021-10-05 14:21:44.711340",
    "2021-10-05 14:21:44.711343",
    "2021-10-05 14:21:44.711345",
    "2021-10-05 14:21:44.711348",
    "2021-10-05 14:21:44.711350",
    "2021-10-05 14:21:44.711353",
    "2021-10-05 14:21:44.711355"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ---------------------------------------------------------------------------
# 2) Add a new "metadata" tab (after "data") describing the panel query that
#    produced this export.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Reuse the same bold / bordered header style already used on "data"!B1:F1
# (format-only copy so no new style entries get minted).
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row
$meta.Range("A2").Value = 0
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$meta.Range("B2").Value = "Neuroendocrine cancer pertinent cancer susceptibility"
$meta.Range("C2").Value = 183

# "1.2" must stay a literal text value (not be auto-converted to the number
# 1.2), but must NOT pick up a new "@" text style either - so build it via a
# text formula in a scratch cell, then paste-special just the computed value
# across (which carries the text type without touching formatting/styles).
$scratch = $meta.Cells.Item(100, 100)
$scratch.Formula = '="1."&"2"'
$scratch.Copy()
$meta.Range("D2").PasteSpecial(-4163)
$scratch.ClearContents()
$excel.CutCopyMode = $false

$meta.Range("E2").Value = "2021-07-28T13:53:35.630017Z"
$meta.Range("F2").Value = "2021-10-05 14:21:44.707750"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/183/?format=json"
